$d = $word.ActiveDocument

# --- Paragraph 1: "**FOR IMMEDIATE RELEASE:" -> bold "FOR IMMEDIATE RELEASE:" ---
$p1 = $d.Paragraphs(1).Range
$t1 = "FOR IMMEDIATE RELEASE:"
$p1.Text = $t1
$d.Range($p1.Start, $p1.Start + $t1.Length).Bold = 1

# --- Paragraph 2: bold "DESIRED NWF DAILY NEWS" and the "SUNDAY LIFESTYLE ..." line ---
$found = $d.Content
$found.Find.Execute("DESIRED NWF DAILY NEWS", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$found.Bold = 1

$found2 = $d.Content
$found2.Find.Execute("SUNDAY LIFESTYLE PUBLICATION DATE: May 5, 2019", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$found2.Bold = 1

# --- Paragraph 3: "Linux User Group Meeting**" -> bold "Linux User Group Meeting" ---
$p3 = $d.Paragraphs(3).Range
$t3 = "Linux User Group Meeting"
$p3.Text = $t3
$d.Range($p3.Start, $p3.Start + $t3.Length).Bold = 1
